# Insert a new weekly data row at row 272 (pushing the existing rows
# 272..391 down to 273..392), then populate the new row with the latest
# week's price observation for "Repollo" (Crespo record, Segunda) at
# Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(272).Insert()

$ws.Cells.Item(272, 1).Value = 4
$ws.Cells.Item(272, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(272, 3).Value = "Los Lagos"
$ws.Cells.Item(272, 4).Value = "2022-03-16"
$ws.Cells.Item(272, 5).Value = 10
$ws.Cells.Item(272, 6).Value = 100112006
$ws.Cells.Item(272, 7).Value = "Repollo"
$ws.Cells.Item(272, 8).Value = "Crespo record"
$ws.Cells.Item(272, 9).Value = "Segunda"
$ws.Cells.Item(272, 10).Value = 100
$ws.Cells.Item(272, 11).Value = 1700
$ws.Cells.Item(272, 12).Value = 1700
$ws.Cells.Item(272, 13).Value = 1700
$ws.Cells.Item(272, 14).Value = "`$/unidad"
$ws.Cells.Item(272, 15).Value = "Región Metropolitana"
$ws.Cells.Item(272, 16).Value = 1700
$ws.Cells.Item(272, 17).Value = 1
$ws.Cells.Item(272, 18).Value = "Hortaliza"
